$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "login" row: EN label gains a trailing full-width chevron
$ws.Range("B3").Value = "Login  ＞"

# "recover_password" row: add the missing JA translation (Submit button)
$ws.Range("C11").Value = "送信"

# "cancel" row: EN label gains a leading full-width chevron + spacing
$ws.Range("B12").Value = "＜　Cancel"

# "recover_password_alert_header" row: add the missing EN translation
$ws.Range("B17").Value = "Nice"
